$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new inventory record (row 4) describing the new asset.
$ws.Range("A4").Value = "Tersedia"
$ws.Range("B4").Value = "Komputer"
$ws.Range("C4").Value = "Laptop Acer Aspire 5"
$ws.Range("D4").Value = "27/12/2023"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "-"
$ws.Range("G4").Value = "Metal"
$ws.Range("H4").Value = 2017
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 8500000
$ws.Range("K4").Value = "Kwarcab Banyumas TA 2017"

# The longer "nama" text no longer fits column C at its old best-fit width;
# re-fit it so the new value is fully visible (mirrors Excel auto-resizing
# the column when the cell content grows).
$ws.Columns("C:C").AutoFit() | Out-Null
$ws.Columns("C:C").ColumnWidth = 18.5

# Restore the selection Excel leaves behind after entering the new data.
$ws.Range("I9").Select()
